$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 77, shifting existing rows 77..166 down to 78..167
$ws.Rows.Item(77).Insert()

# Populate the newly inserted row 77 with the new data record
$ws.Cells.Item(77, 1).Value  = 11
$ws.Cells.Item(77, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(77, 3).Value  = "Bíobío"
$ws.Cells.Item(77, 4).Value  = "01/25/2022"
$ws.Cells.Item(77, 5).Value  = 8
$ws.Cells.Item(77, 6).Value  = 100114001
$ws.Cells.Item(77, 7).Value  = "Papa"
$ws.Cells.Item(77, 8).Value  = "Asterix"
$ws.Cells.Item(77, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(77, 10).Value = 220
$ws.Cells.Item(77, 11).Value = 9000
$ws.Cells.Item(77, 12).Value = 9500
$ws.Cells.Item(77, 13).Value = 9227
$ws.Cells.Item(77, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(77, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(77, 16).Value = 369
$ws.Cells.Item(77, 17).Value = 25
$ws.Cells.Item(77, 18).Value = "Hortaliza"
